# Apply worksheet updates described in the commit diff:
# - Add eight new data columns (G:N) with header labels in row 1
# - Update existing values in columns E:F for several rows
# - Populate the new G:N columns for rows 3-7 and 9-13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers for new columns G:N ---
$ws.Range("G1").Value = "Puerto Rico 2010"
$ws.Range("H1").Value = "Puerto Rico 2020"
$ws.Range("I1").Value = "PR-Born US 2010"
$ws.Range("J1").Value = "PR-Born US 2020"
$ws.Range("K1").Value = "DR-Born US 2010"
$ws.Range("L1").Value = "DR-Born US 2020"
$ws.Range("M1").Value = "Cuban-Born US 2010"
$ws.Range("N1").Value = "Cuban-Born US 2020"

# --- Row 3 ---
$ws.Range("E3").Value = 0.4284603559918569
$ws.Range("G3").Value = 0.2111066509793065
$ws.Range("H3").Value = 0.2019493253082225
$ws.Range("I3").Value = 0.1807980049875312
$ws.Range("J3").Value = 0.1798769971689824
$ws.Range("K3").Value = 0.3160047796443382
$ws.Range("L3").Value = 0.3505808386479016
$ws.Range("M3").Value = 0.1055294721670657
$ws.Range("N3").Value = 0.1309123717153491

# --- Row 4 ---
$ws.Range("E4").Value = 0.3163013565833961
$ws.Range("F4").Value = 0.2853990269601914
$ws.Range("G4").Value = 0.2589660184095785
$ws.Range("H4").Value = 0.2536142122124065
$ws.Range("I4").Value = 0.2977806529700853
$ws.Range("J4").Value = 0.266356447886499
$ws.Range("K4").Value = 0.3244394461235678
$ws.Range("L4").Value = 0.2950804863507435
$ws.Range("M4").Value = 0.2633968613332838
$ws.Range("N4").Value = 0.2237585053193489

# --- Row 5 ---
$ws.Range("E5").Value = 0.2023341068846063
$ws.Range("F5").Value = 0.2217133026906432
$ws.Range("G5").Value = 0.3904934597550003
$ws.Range("H5").Value = 0.3933443355014076
$ws.Range("I5").Value = 0.400145034902284
$ws.Range("J5").Value = 0.4276528586769061
$ws.Range("K5").Value = 0.2624727630561608
$ws.Range("L5").Value = 0.2771223645462593
$ws.Range("M5").Value = 0.4526694750125262
$ws.Range("N5").Value = 0.4484117138453442

# --- Row 6 ---
$ws.Range("E6").Value = 0.03308267106661048
$ws.Range("F6").Value = 0.03384673279036637
$ws.Range("G6").Value = 0.1394338708561146
$ws.Range("H6").Value = 0.1510921269779633
$ws.Range("I6").Value = 0.08572188774924613
$ws.Range("J6").Value = 0.08966841300315642
$ws.Range("K6").Value = 0.06955788289871372
$ws.Range("L6").Value = 0.05436216147856408
$ws.Range("M6").Value = 0.1582942082505985
$ws.Range("N6").Value = 0.1703244238938386

# --- Row 7 ---
$ws.Range("E7").Value = 14239
$ws.Range("G7").Value = 4400
$ws.Range("H7").Value = 13137
$ws.Range("I7").Value = 5312
$ws.Range("J7").Value = 6514
$ws.Range("K7").Value = 1704
$ws.Range("L7").Value = 1997
$ws.Range("M7").Value = 4848
$ws.Range("N7").Value = 6318

# --- Row 9 ---
$ws.Range("E9").Value = 0.4205372064118612
$ws.Range("G9").Value = 0.186765305057294
$ws.Range("H9").Value = 0.1879103069095684
$ws.Range("I9").Value = 0.1680805492464685
$ws.Range("J9").Value = 0.1707538771494167
$ws.Range("K9").Value = 0.2532580071574873
$ws.Range("L9").Value = 0.3276508595823238
$ws.Range("M9").Value = 0.1013584659419538
$ws.Range("N9").Value = 0.1273494988419046

# --- Row 10 ---
$ws.Range("E10").Value = 0.3169901958836399
$ws.Range("F10").Value = 0.273823430547374
$ws.Range("G10").Value = 0.2746508988911599
$ws.Range("H10").Value = 0.2666430740193665
$ws.Range("I10").Value = 0.313707219536816
$ws.Range("J10").Value = 0.285545612144497
$ws.Range("K10").Value = 0.3197461117738414
$ws.Range("L10").Value = 0.2647109726548979
$ws.Range("M10").Value = 0.2495291343487224
$ws.Range("N10").Value = 0.2176296964355258

# --- Row 11 ---
$ws.Range("E11").Value = 0.1955635324842761
$ws.Range("F11").Value = 0.2166665002612565
$ws.Range("G11").Value = 0.3993272206265132
$ws.Range("H11").Value = 0.3858115870671262
$ws.Range("I11").Value = 0.381521816499395
$ws.Range("J11").Value = 0.4087522841212576
$ws.Range("K11").Value = 0.3142767111571271
$ws.Range("L11").Value = 0.2928348909657321
$ws.Range("M11").Value = 0.4338341069866081
$ws.Range("N11").Value = 0.4325025109147929

# --- Row 12 ---
$ws.Range("E12").Value = 0.04780417375870085
$ws.Range("F12").Value = 0.05262071880480979
$ws.Range("G12").Value = 0.1392565754250329
$ws.Range("H12").Value = 0.1596350320039389
$ws.Range("I12").Value = 0.09307035764434868
$ws.Range("J12").Value = 0.08886285901700793
$ws.Range("K12").Value = 0.09840419545792163
$ws.Range("L12").Value = 0.07866043613707165
$ws.Range("M12").Value = 0.1915792908095923
$ws.Range("N12").Value = 0.2030048988460041

# --- Row 13 ---
$ws.Range("E13").Value = 12376
$ws.Range("G13").Value = 3466
$ws.Range("H13").Value = 9792
$ws.Range("I13").Value = 4167
$ws.Range("J13").Value = 4544
$ws.Range("K13").Value = 1046
$ws.Range("L13").Value = 1238
$ws.Range("M13").Value = 4018
$ws.Range("N13").Value = 4705
